$wb = $excel.ActiveWorkbook

# --- Wander through a few sheets, mirroring the author's click-path before
# --- landing on BD/EPG to add the sixth bridge domain + EPG rows.

$wsTenant = $wb.Worksheets.Item("TENANT")
$wsTenant.Activate()
$wsTenant.Range("B2").Select()

$wsLinkLevel = $wb.Worksheets.Item("LINK_LEVEL_POLICY")
$wsLinkLevel.Activate()
$wsLinkLevel.Range("B4").Select()

$wsLacp = $wb.Worksheets.Item("LACP_POLICY")
$wsLacp.Activate()
$wsLacp.Range("C2").Select()

$wsVrf = $wb.Worksheets.Item("VRF")
$wsVrf.Activate()
$wsVrf.Range("B2").Select()

# --- BD sheet: add the sixth bridge domain row ---
$wsBd = $wb.Worksheets.Item("BD")
$wsBd.Activate()
$wsBd.Range("A7").Value = "bd"
$wsBd.Range("B7").Value = "mark_sixth_bd_for_subnet"
$wsBd.Range("C7").Value = "This 6th  bridge domain is created by the Terraform ACI provider"
$wsBd.Range("B7").Select()

# --- EPG sheet: add the sixth EPG row, referencing the new BD ---
$wsEpg = $wb.Worksheets.Item("EPG")
$wsEpg.Activate()
$wsEpg.Range("A7").Value = "epg"
$wsEpg.Range("B7").Value = "mark_sixth_epg"
$wsEpg.Range("C7").Value = "mark_sixth_bd_for_subnet"
$wsEpg.Range("D7").Value = "prod_app_profile"
$wsEpg.Range("A7").Select()

$wsSubnet = $wb.Worksheets.Item("SUBNET")
$wsSubnet.Activate()
$wsSubnet.Range("C2").Select()

# --- Finish with EPG as the active sheet/tab ---
$wsEpg.Activate()
$wsEpg.Range("A7").Select()
